$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3033.1428
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("H40").Value = 29413702
$ws.Range("I40").Value = 1728.2
$ws.Range("J40").Value = 71430810
$ws.Range("K40").Value = 1728.2
$ws.Range("L40").Value = 71430810
$ws.Range("M40").Value = -1553.2
$ws.Range("N40").Value = -71431160
$ws.Range("H41").Value = 956.3333
$ws.Range("I41").Value = 1060.4
$ws.Range("J41").Value = 436
$ws.Range("K41").Value = 1060.4
$ws.Range("L41").Value = 436
$ws.Range("M41").Value = -620.4000000000001
$ws.Range("N41").Value = -1316
$ws.Range("H116").Value = 10340.23
$ws.Range("J116").Value = 10753.956
$ws.Range("L116").Value = 10753.956
$ws.Range("N116").Value = -17637.956
$ws.Range("H132").Value = 6297.3335
$ws.Range("I132").Value = 3321.12
$ws.Range("K132").Value = 9963.360000000001
$ws.Range("M132").Value = -7433.360000000001
$ws.Range("N18").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3188.6924
$ws.Range("I32").Value = 2924.4138
$ws.Range("K32").Value = 2924.4138
$ws.Range("M32").Value = -2637.4138
$ws.Range("H55").Value = 16749.75
$ws.Range("H61").Value = 120001920
$ws.Range("I61").Value = 175002880
$ws.Range("K61").Value = 175002880
$ws.Range("M61").Value = -175002668
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("H102").Value = 2737
$ws.Range("I102").Value = 1862
$ws.Range("K102").Value = 1862
$ws.Range("M102").Value = -240
$ws.Range("H132").Value = 13161532
$ws.Range("I132").Value = 3680
$ws.Range("K132").Value = 11040
$ws.Range("M132").Value = -8510
$ws.Range("H136").Value = 120001920
$ws.Range("I136").Value = 175002880
$ws.Range("K136").Value = 525008640
$ws.Range("M136").Value = -525006090
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 689683.0600000001
$ws.Range("I86").Value = 1168532.4
$ws.Range("K86").Value = 1168532.4
$ws.Range("M86").Value = -1167409.4
$ws.Range("H89").Value = 689683.0600000001
$ws.Range("I89").Value = 1168532.4
$ws.Range("K89").Value = 5842662
$ws.Range("M89").Value = -5837046

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1162
$ws.Range("I22").Value = 682.6667
$ws.Range("K22").Value = 682.6667
$ws.Range("M22").Value = -332.6667
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("H51").Value = 21088
$ws.Range("I51").Value = 21088
$ws.Range("K51").Value = 21088
$ws.Range("M51").Value = -20352
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("N59").Value = 0
$ws.Range("H61").Value = 21088
$ws.Range("I61").Value = 21088
$ws.Range("K61").Value = 21088
$ws.Range("M61").Value = -20740
$ws.Range("H62").Value = 333333340
$ws.Range("J62").Value = 333333340
$ws.Range("L62").Value = 333333340
$ws.Range("N62").Value = -333334588
$ws.Range("H65").Value = 333333340
$ws.Range("J65").Value = 333333340
$ws.Range("L65").Value = 1666666700
$ws.Range("N65").Value = -1666672940
$ws.Range("H68").Value = 99999
$ws.Range("J68").Value = 99999
$ws.Range("L68").Value = 99999
$ws.Range("N68").Value = -101497
$ws.Range("H71").Value = 99999
$ws.Range("J71").Value = 99999
$ws.Range("L71").Value = 299997
$ws.Range("N71").Value = -307485
$ws.Range("H74").Value = 99999
$ws.Range("J74").Value = 99999
$ws.Range("L74").Value = 99999
$ws.Range("N74").Value = -101747
$ws.Range("H77").Value = 99999
$ws.Range("J77").Value = 99999
$ws.Range("L77").Value = 299997
$ws.Range("N77").Value = -308733
$ws.Range("H105").Value = 1296.9166
$ws.Range("I105").Value = 1054.5714
$ws.Range("J105").Value = 2993.3333
$ws.Range("K105").Value = 1054.5714
$ws.Range("L105").Value = 2993.3333
$ws.Range("M105").Value = 692.4286
$ws.Range("N105").Value = -6487.3333
$ws.Range("M50").ClearContents()
$ws.Range("L59").ClearContents()
$ws.Range("M59").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1553.125
$ws.Range("J5").Value = 1527.4
$ws.Range("L5").Value = 4582.200000000001
$ws.Range("N5").Value = -4806.200000000001
$ws.Range("H34").Value = 10904.833
$ws.Range("J34").Value = 16207.75
$ws.Range("L34").Value = 48623.25
$ws.Range("N34").Value = -48791.25
$ws.Range("H55").Value = 5256.3076
$ws.Range("J55").Value = 11958
$ws.Range("L55").Value = 35874
$ws.Range("N55").Value = -36228
$ws.Range("H135").Value = 1553.125
$ws.Range("J135").Value = 1527.4
$ws.Range("L135").Value = 13746.6
$ws.Range("N135").Value = -18816.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 9999
$ws.Range("I52").Value = 9999
$ws.Range("K52").Value = 9999
$ws.Range("M52").Value = -9740
$ws.Range("H126").Value = 10843.454
$ws.Range("I126").Value = 13996.875
$ws.Range("K126").Value = 41990.625
$ws.Range("M126").Value = -39520.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 996.3333
$ws.Range("I46").Value = 996.3333
$ws.Range("K46").Value = 996.3333
$ws.Range("M46").Value = -808.3333
$ws.Range("H61").Value = 71431690
$ws.Range("I61").Value = 76926170
$ws.Range("K61").Value = 76926170
$ws.Range("M61").Value = -76925968
$ws.Range("H68").Value = 2453390.5
$ws.Range("I68").Value = 5210831
$ws.Range("J68").Value = 2332.2222
$ws.Range("K68").Value = 5210831
$ws.Range("L68").Value = 2332.2222
$ws.Range("M68").Value = -5210082
$ws.Range("N68").Value = -3830.2222
$ws.Range("H71").Value = 2453390.5
$ws.Range("I71").Value = 5210831
$ws.Range("J71").Value = 2332.2222
$ws.Range("K71").Value = 26054155
$ws.Range("L71").Value = 11661.111
$ws.Range("M71").Value = -26050411
$ws.Range("N71").Value = -19149.111
$ws.Range("H113").Value = 71431690
$ws.Range("I113").Value = 76926170
$ws.Range("K113").Value = 76926170
$ws.Range("M113").Value = -76924000
$ws.Range("H136").Value = 2548.94
$ws.Range("I136").Value = 2326.9697
$ws.Range("J136").Value = 2979.8235
$ws.Range("K136").Value = 6980.909100000001
$ws.Range("L136").Value = 8939.470499999999
$ws.Range("M136").Value = -4430.909100000001
$ws.Range("N136").Value = -14039.4705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10980.6
$ws.Range("J41").Value = 10980.6
$ws.Range("L41").Value = 10980.6
$ws.Range("N41").Value = -11760.6
$ws.Range("H96").Value = 16103.637
$ws.Range("I96").Value = 9073.333000000001
$ws.Range("J96").Value = 24540
$ws.Range("K96").Value = 9073.333000000001
$ws.Range("L96").Value = 24540
$ws.Range("M96").Value = -7700.333000000001
$ws.Range("N96").Value = -27286
$ws.Range("H136").Value = 7608.718
$ws.Range("I136").Value = 7760.2646
$ws.Range("K136").Value = 23280.7938
$ws.Range("M136").Value = -20730.7938
